$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to force Excel to store numeric-looking strings as text
# (matching the source inlineStr cells) instead of auto-converting them to numbers.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

function Set-TextValue($cellAddr, $text) {
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
}

Set-TextValue "D2" '30.164.74'
$ws.Range("E2").Value = '  -0.94%  '

Set-TextValue "D3" '1.849.75'
$ws.Range("E3").Value = '  -2.24%  '

Set-TextValue "D4" '1.000'
$ws.Range("E4").Value = '  -0.18%  '

Set-TextValue "D5" '236.28'
$ws.Range("E5").Value = '  -0.68%  '

Set-TextValue "D6" '0.9998'
$ws.Range("E6").Value = '  -0.16%  '

Set-TextValue "D7" '0.4762'
$ws.Range("E7").Value = '  -2.73%  '

Set-TextValue "D8" '0.2821'
$ws.Range("E8").Value = '  -3.79%  '

Set-TextValue "D9" '0.06481'
$ws.Range("E9").Value = '  -2.98%  '

Set-TextValue "D10" '1.844.75'
$ws.Range("E10").Value = '  -2.44%  '

$ws.Range("E11").Value = '  -0.70%  '

Set-TextValue "D12" '16.35'
$ws.Range("E12").Value = '  -4.22%  '

Set-TextValue "D13" '5.133'
$ws.Range("E13").Value = '  +0.00%  '

Set-TextValue "D14" '87.24'
$ws.Range("E14").Value = '  -0.87%  '

Set-TextValue "D15" '0.6455'
$ws.Range("E15").Value = '  -2.71%  '

Set-TextValue "D16" '30.116.49'
$ws.Range("E16").Value = '  -1.00%  '

$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue "D17" '13.23'
$ws.Range("E17").Value = '  -1.51%  '

$ws.Range("B18").Value = 'Dai'
$ws.Range("C18").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue "D18" '0.9996'
$ws.Range("E18").Value = '  -0.20%  '

Set-TextValue "D19" '0.000007629'
$ws.Range("E19").Value = '  -2.44%  '

Set-TextValue "D20" '2.107.39'
$ws.Range("E20").Value = '  -1.07%  '

Set-TextValue "D21" '1.001'
$ws.Range("E21").Value = '  -0.02%  '

Set-TextValue "D22" '5.274'
$ws.Range("E22").Value = '  -0.17%  '

Set-TextValue "D23" '217.85'
$ws.Range("E23").Value = '  +15.93%  '

Set-TextValue "D24" '6.101'
$ws.Range("E24").Value = '  -0.75%  '

Set-TextValue "D25" '9.189'
$ws.Range("E25").Value = '  -2.99%  '

Set-TextValue "D26" '163.77'
$ws.Range("E26").Value = '  -0.04%  '

Set-TextValue "D27" '18.38'
$ws.Range("E27").Value = '  +0.46%  '

Set-TextValue "D28" '1.913'
$ws.Range("E28").Value = '  -0.76%  '

Set-TextValue "D29" '1.427'
$ws.Range("E29").Value = '  -2.62%  '

Set-TextValue "D30" '0.09186'
$ws.Range("E30").Value = '  +0.16%  '

Set-TextValue "D31" '4.240'
$ws.Range("E31").Value = '  -2.40%  '

Set-TextValue "D32" '3.963'
$ws.Range("E32").Value = '  -2.93%  '

Set-TextValue "D33" '0.05012'
$ws.Range("E33").Value = '  -3.61%  '

Set-TextValue "D34" '0.7412'
$ws.Range("E34").Value = '  +0.17%  '

Set-TextValue "D35" '1.137'
$ws.Range("E35").Value = '  +3.64%  '

Set-TextValue "D36" '2.685'
$ws.Range("E36").Value = '  -1.19%  '

Set-TextValue "D37" '0.01830'
$ws.Range("E37").Value = '  +0.72%  '

Set-TextValue "D38" '2.609'
$ws.Range("E38").Value = '  -2.37%  '

Set-TextValue "D39" '2.068'
$ws.Range("E39").Value = '  +1.65%  '

Set-TextValue "D40" '0.9017'
$ws.Range("E40").Value = '  -1.58%  '

Set-TextValue "D41" '5.917'
$ws.Range("E41").Value = '  -0.50%  '

Set-TextValue "D42" '106.27'
$ws.Range("E42").Value = '  +0.61%  '

$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue "D43" '0.4245'
$ws.Range("E43").Value = '  -3.35%  '

$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue "D44" '0.9987'
$ws.Range("E44").Value = '  +0.58%  '

Set-TextValue "D45" '7.426'
$ws.Range("E45").Value = '  -1.93%  '

Set-TextValue "D46" '0.1306'
$ws.Range("E46").Value = '  -5.65%  '

Set-TextValue "D47" '1.559'
$ws.Range("E47").Value = '  +10.59%  '

Set-TextValue "D48" '63.88'
$ws.Range("E48").Value = '  -6.06%  '

Set-TextValue "D49" '8.774'
$ws.Range("E49").Value = '  -2.75%  '

Set-TextValue "D50" '34.23'
$ws.Range("E50").Value = '  -2.02%  '

$ws.Range("E51").Value = '  -2.56%  '

$scratch.Clear()